# "Able to continue or start new game, also round infection percent"
#
# This inserts a new "Sprint 6" worksheet (burndown chart + backlog items
# that moved from "in options/planning" into an active sprint), wires the
# Product Backlog roll-up to it, and marks the two moved items "In Progress"
# with a yellow status fill.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert "Sprint 6" right after "Product Backlog" (before "Sprint 5")
# ---------------------------------------------------------------------
# NOTE: worksheet variables in this host are index-bound, not identity-
# bound — once a sheet is inserted/removed the indices shift underneath
# any previously-held reference. So every sheet handle below is re-
# fetched by name (via Worksheets.Item) right before it's used, rather
# than reused across the structural Add() call.
$sprint5 = $wb.Worksheets.Item("Sprint 5")
$sprint6 = $wb.Worksheets.Add($sprint5)
$sprint6.Name = "Sprint 6"
$sprint6 = $wb.Worksheets.Item("Sprint 6")

# ---------------------------------------------------------------------
# 2. Populate "Sprint 6" — burndown data + the two carried-over items
# ---------------------------------------------------------------------

# Column widths (approximate the source's auto-fit widths)
$sprint6.Columns("A:A").ColumnWidth = 25.27
$sprint6.Columns("B:L").ColumnWidth = 10.45
$sprint6.Columns("M:P").ColumnWidth = 9.45

# Row 1: sprint date headers, B1:O1 = 44156..44169
$startDate = 44156
for ($i = 0; $i -lt 14; $i++) {
    $cell = $sprint6.Cells.Item(1, 2 + $i)
    $cell.Value = $startDate + $i
    $cell.NumberFormat = "m/d/yy"
}
$sprint6.Cells.Item(1, 16).NumberFormat = "m/d/yy"

# Row 2: Optimal Trend
$sprint6.Range("A2").Value = "Optimal Trend"
$sprint6.Range("B2").Value = 8
$sprint6.Range("C2").Formula = "=B2 - `$B`$5"
$sprint6.Range("C2").NumberFormat = "0.00"
$prevCol = "C"
$cols = @("D","E","F","G","H","I","J","K","L","M","N","O")
foreach ($col in $cols) {
    $f = "=" + $prevCol + "2 - `$B`$5"
    $sprint6.Range($col + "2").Formula = $f
    $sprint6.Range($col + "2").NumberFormat = "0.00"
    $prevCol = $col
}

# Row 3: Actual Trend
$sprint6.Range("A3").Value = "Actual Trend"
$actual = @(8,8,8,8,6,6,6,6,6,6,5,5,3,0)
for ($i = 0; $i -lt $actual.Length; $i++) {
    $sprint6.Cells.Item(3, 2 + $i).Value = $actual[$i]
}

# Row 5: Hours Per Day
$sprint6.Range("A5").Value = "Hours Per Day"
$sprint6.Range("B5").Formula = "=8/13"
$sprint6.Range("B5").NumberFormat = "0.00"

# Row 7: table header
$sprint6.Range("A7").Value = "Use Cases"
$sprint6.Range("B7").Value = "Estimate"
$sprint6.Range("C7").Value = "Status"

# Row 8 & 9: items carried into this sprint, marked "In Progress" (yellow)
$sprint6.Range("A8").Value = "Create options menu"
$sprint6.Range("B8").Value = 4
$sprint6.Range("C8").Value = "In Progress"
$sprint6.Range("C8").Interior.Color = 65535

$sprint6.Range("A9").Value = "Resize start screen for web gl"
$sprint6.Range("B9").Value = 4
$sprint6.Range("C9").Value = "In Progress"
$sprint6.Range("C9").Interior.Color = 65535

# ---------------------------------------------------------------------
# 3. "Sprint 5" tab view: no longer the active tab, new selection A7:C9
# ---------------------------------------------------------------------
$sprint5 = $wb.Worksheets.Item("Sprint 5")
$sprint5.Activate()
$sprint5.Range("A7:C9").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. "Product Backlog": drop the items that moved into Sprint 6 / were
#    otherwise removed, and roll the "Create options menu" & "Resize
#    start screen for web gl" estimate total up from the new Sprint 6 sheet
# ---------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Product Backlog")
$backlog.Range("A3:B4").ClearContents() | Out-Null
$backlog.Range("A6:B8").ClearContents() | Out-Null
$backlog.Range("A10:B10").ClearContents() | Out-Null
$backlog.Range("B5").Formula = "=SUM('Sprint 6'!B8:B9)"
$backlog.Activate()
$backlog.Range("A2:B3").Select() | Out-Null

# ---------------------------------------------------------------------
# 5. View: make Sprint 6 the active/selected tab, scrolled to column B,
#    with M5 selected — done last so it "wins" as the saved active tab.
# ---------------------------------------------------------------------
$sprint6 = $wb.Worksheets.Item("Sprint 6")
$sprint6.Activate()
$sprint6.Application.ActiveWindow.ScrollColumn = 2
$sprint6.Range("M5").Select() | Out-Null

Write-Output "Sprint 6 inserted and wired up"
